{"js": "// Add a centered, bold Arial 12pt title paragraph\n// \"PROYECTO FINAL DE ANALISIS DE SISTEMAS\" as the very first\n// paragraph of the document body (commit: \"se agrego titulo al informe\").\n\nconst body = context.document.body;\n\n// Insert a brand-new paragraph before the current first paragraph.\nconst titlePara = body.insertParagraph(\n  \"PROYECTO FINAL DE ANALISIS DE SISTEMAS\",\n  Word.InsertLocation.start\n);\n\n// Center it and make the text Bold Arial 12pt, matching the rest of the\n// document's heading style.\ntitlePara.alignment = Word.Alignment.centered;\ntitlePara.font.name = \"Arial\";\ntitlePara.font.bold = true;\ntitlePara.font.size = 12;\n\nawait context.sync();\n", "ps1": "# Add a centered, bold Arial 12pt title paragraph\n# \"PROYECTO FINAL DE ANALISIS DE SISTEMAS\" as the very first\n# paragraph of the document body (commit: \"se agrego titulo al informe\").\n\n$d = $word.ActiveDocument\n\n# Zero-length range at the very start of the story; InsertBefore there\n# creates a brand-new first paragraph (the trailing carriage return ends it).\n$r = $d.Range(0, 0)\n$r.InsertBefore(\"PROYECTO FINAL DE ANALISIS DE SISTEMAS`r\")\n\n# Format the newly created first paragraph: centered, Bold Arial 12pt.\n$p = $d.Paragraphs(1)\n$p.Alignment = 1\n$p.Range.Font.Name = \"Arial\"\n$p.Range.Font.Bold = $true\n$p.Range.Font.Size = 12\n"}
